$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Helper: split the run that currently occupies [start, start+len) into
# two runs at offset `splitOffset` (relative to start) by temporarily
# dropping a bookmark at the split point and immediately removing it.
# Word (and this COM host) never re-merges runs that were once split by
# a bookmark, even after the bookmark itself is deleted.
# ---------------------------------------------------------------------
function Split-RunAt($absolutePos, $bookmarkName) {
    $bmRange = $d.Range($absolutePos, $absolutePos)
    $d.Bookmarks.Add($bookmarkName, $bmRange)
    $bm = $d.Bookmarks.Item($bookmarkName)
    $bm.Delete()
}

# ------------------------------------------------------------------
# 1) First bullet ("la stanza non può essere un balcone ..." block):
#    "potrebbe essere richiesto oltre al pavimento anche il rivestimento"
#    becomes two runs:
#      "Chiedere se fare oltre " + "al pavimento anche il rivestimento"
# ------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("potrebbe essere richiesto oltre al pavimento anche il rivestimento", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Chiedere se fare oltre al pavimento anche il rivestimento", 1) | Out-Null

$rng2 = $d.Content
$rng2.Find.Execute("Chiedere se fare oltre al pavimento anche il rivestimento", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitPos = $rng2.Start + "Chiedere se fare oltre ".Length
Split-RunAt $splitPos "TmpSplitA"

# ------------------------------------------------------------------
# 2) Second bullet ("la stanza sarà un balcone o terrazzo ..." block):
#    "potrebbe essere richiesto oltre al pavimento anche il rivestimento"
#    + " (in particolare)" becomes two runs:
#      "solo " + "pavimento"
# ------------------------------------------------------------------
$rng3 = $d.Content
$rng3.Find.Execute("potrebbe essere richiesto oltre al pavimento anche il rivestimento (in particolare)", `
    $true, $false, $false, $false, $false, $true, 1, $false, "solo pavimento", 1) | Out-Null

$rng4 = $d.Content
$rng4.Find.Execute("solo pavimento", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitPos2 = $rng4.Start + "solo ".Length
Split-RunAt $splitPos2 "TmpSplitB"

# ------------------------------------------------------------------
# 3) Remove the stray _GoBack bookmark that currently sits, on its own,
#    in the otherwise-empty paragraph right after "Mostrare foto dei
#    vari tipi di posa".
# ------------------------------------------------------------------
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# ------------------------------------------------------------------
# 4) Re-create the _GoBack bookmark inside "Ci si troverà a realizzare
#    una forma triangolare. ", splitting it into "Ci si tr" | bookmark |
#    "overà a realizzare una forma triangolare. "
# ------------------------------------------------------------------
$rng5 = $d.Content
$rng5.Find.Execute("Ci si troverà a realizzare una forma triangolare.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitPos3 = $rng5.Start + "Ci si tr".Length
$bmRange3 = $d.Range($splitPos3, $splitPos3)
$d.Bookmarks.Add("_GoBack", $bmRange3)
